$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.5440079918808879; C = 0.5440079918808879; D = 0.3686258812966464; E = 0.6071456837503223; F = 0.2797714183633027; G = 14 }
    3  = @{ B = 0.371620442828777;  C = 0.3874993246909071; D = 0.1981274257140365; E = 0.4451150701942549; F = 0.2550055127533425; G = 13 }
    4  = @{ B = 0.3184154172361096; C = 0.3634107558758037; D = 0.1745848588293675; E = 0.4178335300444035; F = 0.2825786729951437; G = 12 }
    5  = @{ B = 0.3975283727593981; C = 0.4145492339343189; D = 0.2091250642104284; E = 0.4573019398717093; F = 0.2370777989771771; G = 11 }
    6  = @{ B = 0.3560805627268342; C = 0.3805953001122874; D = 0.1781696942431756; E = 0.4221015212519088; F = 0.2389242722689363; G = 10 }
    7  = @{ B = 0.3392710727664942; C = 0.3699241520168967; D = 0.1693631219851699; E = 0.4115375098155329; F = 0.2470638456253208; G = 9 }
    8  = @{ B = 0.352534700000344;  C = 0.3776340216059933; D = 0.1795390275838225; E = 0.4237204592462139; F = 0.251301328242729;  G = 8 }
    9  = @{ B = 0.323074747095869;  C = 0.3439312876808158; D = 0.1492574881197667; E = 0.3863385666999435; F = 0.2288235751843726; G = 7 }
    10 = @{ B = 0.3581453753460085; C = 0.3775233390900978; D = 0.1783793918318441; E = 0.4223498453081807; F = 0.2452214067738235; G = 6 }
    11 = @{ B = 0.3599616089605324; C = 0.3766256120414225; D = 0.1882494404226522; E = 0.4338772181420132; F = 0.2708253138491602; G = 5 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
